# SU suite assemblage Front Rocker
# mise en place des vis pour BAR, Triangle et Damper.
# Rectification : suppression des rondelles côté écrou Knut

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# rayon_vis_basculeurs_suspension (mm) : 4 -> 4.2
$ws.Range("B23").Value = 4.2

# rayon_vis_basculeurs_autres (mm) : 3 -> 4.2
$ws.Range("B24").Value = 4.2

# Update the active cell selection to B24
$ws.Range("B24").Select()
